$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Transmitance column (B3:B18) from 1 to 100 (extra ordinary ray for polarimetric acquisitions)
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 2).Value = 100
}

# Update the active cell selection to B18
$ws.Range("B18").Select()
